$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New "Quinte Flush" rule row (row 21) has now been tested: copy the date
# formatting used by the row above (H20:I20) onto H21:I21, then fill in the
# test start/end dates and the "OK" result.
$ws.Range("H20:I20").Copy()
$ws.Range("H21:I21").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("H21").Value = 43173
$ws.Range("I21").Value = 43173
$ws.Range("J21").Value = "OK"

# Reflect where the user ended up looking in the sheet after the edit.
$ws.Range("I26").Select()
